$d = $word.ActiveDocument

$replacements = @(
    @{old="728×9="; new="677×9="},
    @{old="806×7="; new="893×6="},
    @{old="295×9="; new="168×5="},
    @{old="228×8="; new="407×3="},
    @{old="734×2="; new="107×9="},
    @{old="407×6="; new="712×6="},
    @{old="149×6="; new="101×4="},
    @{old="895×8="; new="411×3="},
    @{old="469×3="; new="587×8="},
    @{old="812×6="; new="722×6="},
    @{old="467×5="; new="644×8="},
    @{old="675×8="; new="430×6="},
    @{old="712×7="; new="292×4="},
    @{old="621×6="; new="823×7="},
    @{old="840×3="; new="660×6="},
    @{old="579×4="; new="367×3="},
    @{old="688×5="; new="262×8="},
    @{old="968×7="; new="735×3="},
    @{old="119×2="; new="342×8="},
    @{old="238×7="; new="244×2="},
    @{old="120×6="; new="854×9="},
    @{old="944×2="; new="770×8="},
    @{old="651×5="; new="401×7="},
    @{old="518×5="; new="939×3="},
    @{old="172×2="; new="664×9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
